# fixing relation for rev 0
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("components request")

# Duplicate the existing "test_rev_A2" row down to row 3 before touching row 2,
# so the old value/formula pattern is preserved there.
$ws.Range("B3").Value = "test_rev_A2"
$ws.Range("C3").Formula = "=RIGHT(B3,IF(RIGHT(B3,1)=""0"",,2))"

# Row 2 becomes the new "rev 0" test case and the formula drops the erroneous
# literal 0 in favor of an empty (omitted) argument.
$ws.Range("B2").Value = "test_rev_0"
$ws.Range("C2").Formula = "=RIGHT(B2,IF(RIGHT(B2,1)=""0"",,2))"

$ws.Range("B2").Select()
